$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The poll-answer sheet had a "Question"/"Answer" header row in row 2
# (the real header with the sheet title is row 1). This update removes
# that redundant header row, shifting every question/answer pair up by
# one row and shrinking the used range from A1:B12 to A1:B11.
$ws.Rows.Item(2).Delete()

# Reflect the selection left behind by the edit (row 2, now the first
# question, is selected as a full row).
$ws.Range("A2:XFD2").Select()
